$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two "外籍學者..." header labels to "境外學者..."
$ws.Range("C1").Value = "境外學者姓名"
$ws.Range("D1").Value = "境外學者身分（教授、副教授、助理教授或博士後研究員）"

# Update the sheet view: zoom to 100 (normal) and move selection to D1
$excel.ActiveWindow.Zoom = 100
$ws.Range("D1").Select()
